$d = $word.ActiveDocument

# Locate the paragraph that contains the Eid date sentence with the typo
# ("Sunday" should read "Saturday").
$rng = $d.Content
[void]$rng.Find.Execute("This year, Eid will either take place on")
$para = $rng.Paragraphs(1)
$prng = $para.Range

# Rebuild the paragraph's OOXML, splitting the run that held
# " Friday 20 March or Sunday 21 March 2026" into three runs so that
# only the word "Sunday" changes to "Saturday" (matching how Word
# isolates a corrected word into its own run after a typo fix).
$newParaXml = '<w:p w14:paraId="5FC330E9" w14:textId="77777777" w:rsidR="00B43858" w:rsidRDefault="00C75979" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:after="200"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>This year, Eid will either take place on</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="white"/></w:rPr><w:t xml:space="preserve"> Friday 20 March or </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="white"/></w:rPr><w:t>Saturday</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="white"/></w:rPr><w:t xml:space="preserve"> 21 March 2026</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> depending on the sighting of the moon.</w:t></w:r></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$prng.InsertXML($xml)
